$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "...eine INVITE-Nachricht  durch unsere Partnergruppe (IP ...) an..."
#    used to be split across two runs ("-Nachricht  " / "durch unsere
#    Partnergruppe (IP 141.22.27.35) "); the edit just joins that text
#    back into a single run (no wording change).
# ---------------------------------------------------------------------
$merge = $d.Content
$merge.Find.Execute(
    "-Nachricht  durch unsere Partnergruppe (IP 141.22.27.35) ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-Nachricht  durch unsere Partnergruppe (IP 141.22.27.35) ",
    2)

# ---------------------------------------------------------------------
# 2) Extend the closing sentence of the "Fazit" section:
#    "... Somit werden die Multicastpakete lediglich gebroadcastet."
#    becomes
#    "... Somit werden die Multicastpakete lediglich gebroadcastet und
#     erreichen daher auch nach dem IGMP-Leave weiterhin den UAC."
#    ("IGMP-Leave" rendered bold, like the other IGMP-Join/IGMP-Leave
#    mentions throughout the document).
# ---------------------------------------------------------------------

# Remember a stable anchor (before any edits) so later Find calls can be
# restricted to the tail of the document and not accidentally match one
# of the earlier "IGMP-Leave" / "IGMP-Join" occurrences.
$anchor = $d.Content
$anchor.Find.Execute("Somit werden die Multicastpakete lediglich")
$anchorStart = $anchor.Start

# Insert the new (still unformatted) text right before the full stop
# that follows "gebroadcastet".
$tail = $d.Range($anchorStart, $d.Content.End)
$tail.Find.Execute("gebroadcastet.")
$periodPos = $tail.End - 1
$insertionPoint = $d.Range($periodPos, $periodPos)
$insertionPoint.InsertBefore(" und erreichen daher auch nach dem IGMP-Leave weiterhin den UAC")

# Bold only the newly-added "IGMP-" / "Leave" pieces, restricted to the
# tail of the document so the pre-existing bold "IGMP-Join"/"IGMP-Leave"
# occurrences earlier in the text stay untouched.
$boldIgmp = $d.Range($anchorStart, $d.Content.End)
$boldIgmp.Find.Execute("IGMP-")
$boldIgmp.Font.Bold = 1

$boldLeave = $d.Range($anchorStart, $d.Content.End)
$boldLeave.Find.Execute("Leave")
$boldLeave.Font.Bold = 1
